# Xoa lan chon truoc (Remove previous selection)
# Clears the previously-filled "references" text that had been entered in
# F25 and G25 ("3.11 Further Reading" / the Chung-Nixon-Yu-Mylopoulos
# citation), leaving those two cells blank again while keeping their
# formatting (fill/border/alignment) untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F25:G25").ClearContents()

# Leave the selection where the user ended up after clearing the cells.
$ws.Range("F30").Select()
